# Generate Report for Handoff
# Updates the localization status report:
#  - Status moves from "In Translation" to "Ready for handoff"
#  - Refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
#  - Column widths for the Status column(s) widen to fit the new, longer status text

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# E2 = zh-cn status, F2 = de-de status
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
# G2 = Latest HO Xliff Generate Date
$ws_overview.Range("G2").Value = "2016-08-28 11:07:05"

# Widen the two status columns (E, F) to fit the longer text
$ws_overview.Columns.Item(5).ColumnWidth = 16.35
$ws_overview.Columns.Item(6).ColumnWidth = 16.35

# --- zh-cn sheet --------------------------------------------------------
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_zhcn.Range("H2").Value = "2016-08-28 11:06:58"
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.35

# --- de-de sheet --------------------------------------------------------
$ws_dede.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("H2").Value = "2016-08-28 11:07:05"
$ws_dede.Columns.Item(3).ColumnWidth = 16.35
